$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 12.84399113027431
$ws.Cells.Item(2, 3).Value = 8.8400742040184
$ws.Cells.Item(2, 5).Value = 21.42273286564289
$ws.Cells.Item(2, 6).Value = 38.0192710486092
$ws.Cells.Item(2, 7).Value = 3.617520742673954
$ws.Cells.Item(2, 10).Value = 7.705357923337286
$ws.Cells.Item(2, 13).Value = 19.28082668016616
$ws.Cells.Item(2, 15).Value = 19.85006442099378

$ws.Cells.Item(3, 2).Value = 12.21591784434323
$ws.Cells.Item(3, 3).Value = 8.414586927430291
$ws.Cells.Item(3, 5).Value = 21.38088291392584
$ws.Cells.Item(3, 6).Value = 38.0197488058127
$ws.Cells.Item(3, 7).Value = 3.619657161654706
$ws.Cells.Item(3, 10).Value = 7.728478295074376
$ws.Cells.Item(3, 13).Value = 19.02889191414526
$ws.Cells.Item(3, 15).Value = 19.97235111230211

$ws.Cells.Item(4, 2).Value = 11.81370258593417
$ws.Cells.Item(4, 3).Value = 8.1409774274368
$ws.Cells.Item(4, 5).Value = 21.35889497832381
$ws.Cells.Item(4, 6).Value = 38.03171301837283
$ws.Cells.Item(4, 7).Value = 3.621037052954156
$ws.Cells.Item(4, 10).Value = 7.743398818283126
$ws.Cells.Item(4, 13).Value = 18.87491905588814
$ws.Cells.Item(4, 15).Value = 20.05410230866446

$ws.Cells.Item(5, 2).Value = 11.64581255457273
$ws.Cells.Item(5, 3).Value = 8.026467854378669
$ws.Cells.Item(5, 5).Value = 21.35087271591757
$ws.Cells.Item(5, 6).Value = 38.03951865954831
$ws.Cells.Item(5, 7).Value = 3.621616555530624
$ws.Cells.Item(5, 10).Value = 7.749661816928121
$ws.Cells.Item(5, 13).Value = 18.81241841340086
$ws.Cells.Item(5, 15).Value = 20.08908454575149

$ws.Cells.Item(6, 2).Value = 11.61769990055267
$ws.Cells.Item(6, 3).Value = 8.007274882054295
$ws.Cells.Item(6, 5).Value = 21.34959742688777
$ws.Cells.Item(6, 6).Value = 38.04099157438479
$ws.Cells.Item(6, 7).Value = 3.621713821110813
$ws.Cells.Item(6, 10).Value = 7.750712838582016
$ws.Cells.Item(6, 13).Value = 18.80205687987184
$ws.Cells.Item(6, 15).Value = 20.09499381444663

$ws.Cells.Item(7, 2).Value = 11.8114542278342
$ws.Cells.Item(7, 3).Value = 8.139445165314164
$ws.Cells.Item(7, 5).Value = 21.35878298255729
$ws.Cells.Item(7, 6).Value = 38.0318064318197
$ws.Cells.Item(7, 7).Value = 3.621044798670437
$ws.Cells.Item(7, 10).Value = 7.743482542506955
$ws.Cells.Item(7, 13).Value = 18.87407507517312
$ws.Cells.Item(7, 15).Value = 20.054567349002

$ws.Cells.Item(8, 2).Value = 12.63097300065035
$ws.Cells.Item(8, 3).Value = 8.695983933137679
$ws.Cells.Item(8, 5).Value = 21.40753693136047
$ws.Cells.Item(8, 6).Value = 38.01701138500127
$ws.Cells.Item(8, 7).Value = 3.618243275118851
$ws.Cells.Item(8, 10).Value = 7.71317980150515
$ws.Cells.Item(8, 13).Value = 19.19384964809366
$ws.Cells.Item(8, 15).Value = 19.89084013167057

$ws.Cells.Item(9, 2).Value = 14.09987393193757
$ws.Cells.Item(9, 3).Value = 9.685848609206964
$ws.Cells.Item(9, 5).Value = 21.53229234208833
$ws.Cells.Item(9, 6).Value = 38.08070824077436
$ws.Cells.Item(9, 7).Value = 3.61328744731133
$ws.Cells.Item(9, 10).Value = 7.65947850385227
$ws.Cells.Item(9, 13).Value = 19.82378162338566
$ws.Cells.Item(9, 15).Value = 19.62304601905533

$ws.Cells.Item(10, 2).Value = 15.08764655886571
$ws.Cells.Item(10, 3).Value = 10.34753174075186
$ws.Cells.Item(10, 5).Value = 21.64131332290061
$ws.Cells.Item(10, 6).Value = 38.18402558234128
$ws.Cells.Item(10, 7).Value = 3.609970755088169
$ws.Cells.Item(10, 10).Value = 7.623475403984405
$ws.Cells.Item(10, 13).Value = 20.28458345794194
$ws.Cells.Item(10, 15).Value = 19.45929225882637

$ws.Cells.Item(11, 2).Value = 15.51604702430854
$ws.Cells.Item(11, 3).Value = 10.63375053888137
$ws.Cells.Item(11, 5).Value = 21.6945790634793
$ws.Cells.Item(11, 6).Value = 38.243250128863
$ws.Cells.Item(11, 7).Value = 3.60853157016232
$ws.Cells.Item(11, 10).Value = 7.607838142840871
$ws.Cells.Item(11, 13).Value = 20.49300172091274
$ws.Cells.Item(11, 15).Value = 19.39207209634869

$ws.Cells.Item(12, 2).Value = 15.67518858977731
$ws.Cells.Item(12, 3).Value = 10.73997274715429
$ws.Cells.Item(12, 5).Value = 21.71526768156609
$ws.Cells.Item(12, 6).Value = 38.26742654714906
$ws.Cells.Item(12, 7).Value = 3.607996537865269
$ws.Cells.Item(12, 10).Value = 7.602022638283527
$ws.Cells.Item(12, 13).Value = 20.57168836237427
$ws.Cells.Item(12, 15).Value = 19.36767255139902

$ws.Cells.Item(13, 2).Value = 15.64105280316365
$ws.Cells.Item(13, 3).Value = 10.71719258611608
$ws.Cells.Item(13, 5).Value = 21.71078914254517
$ws.Cells.Item(13, 6).Value = 38.26214208197305
$ws.Cells.Item(13, 7).Value = 3.608111324678183
$ws.Cells.Item(13, 10).Value = 7.603270405826419
$ws.Cells.Item(13, 13).Value = 20.55475331588669
$ws.Cells.Item(13, 15).Value = 19.37288034493872

$ws.Cells.Item(14, 2).Value = 15.52920188201742
$ws.Cells.Item(14, 3).Value = 10.64253303098575
$ws.Cells.Item(14, 5).Value = 21.69627080544551
$ws.Cells.Item(14, 6).Value = 38.24520413526603
$ws.Cells.Item(14, 7).Value = 3.608487353499646
$ws.Cells.Item(14, 10).Value = 7.607357576662159
$ws.Cells.Item(14, 13).Value = 20.49948042374413
$ws.Cells.Item(14, 15).Value = 19.39004352737059

$ws.Cells.Item(15, 2).Value = 15.46028635269215
$ws.Cells.Item(15, 3).Value = 10.59651929286911
$ws.Cells.Item(15, 5).Value = 21.68744506293399
$ws.Cells.Item(15, 6).Value = 38.23505667670189
$ws.Cells.Item(15, 7).Value = 3.60871897701625
$ws.Cells.Item(15, 10).Value = 7.609874873801263
$ws.Cells.Item(15, 13).Value = 20.46559147584804
$ws.Cells.Item(15, 15).Value = 19.40069418222171

$ws.Cells.Item(16, 2).Value = 15.05922161816931
$ws.Cells.Item(16, 3).Value = 10.32852586871044
$ws.Cells.Item(16, 5).Value = 21.63790532887112
$ws.Cells.Item(16, 6).Value = 38.18040044069376
$ws.Cells.Item(16, 7).Value = 3.610066205194774
$ws.Cells.Item(16, 10).Value = 7.624512195814943
$ws.Cells.Item(16, 13).Value = 20.27093312836951
$ws.Cells.Item(16, 15).Value = 19.46383240593909

$ws.Cells.Item(17, 2).Value = 14.8077612090223
$ws.Cells.Item(17, 3).Value = 10.16030551300714
$ws.Cells.Item(17, 5).Value = 21.60844782494135
$ws.Cells.Item(17, 6).Value = 38.14999640887998
$ws.Cells.Item(17, 7).Value = 3.610910473971301
$ws.Cells.Item(17, 10).Value = 7.633681059504088
$ws.Cells.Item(17, 13).Value = 20.15116244320573
$ws.Cells.Item(17, 15).Value = 19.50443504516902

$ws.Cells.Item(18, 2).Value = 14.66116004275124
$ws.Cells.Item(18, 3).Value = 10.06215973799766
$ws.Cells.Item(18, 5).Value = 21.59185086782385
$ws.Cells.Item(18, 6).Value = 38.13366015645315
$ws.Cells.Item(18, 7).Value = 3.611402628730461
$ws.Cells.Item(18, 10).Value = 7.639024499693478
$ws.Cells.Item(18, 13).Value = 20.08216401053783
$ws.Cells.Item(18, 15).Value = 19.5284725379831

$ws.Cells.Item(19, 2).Value = 14.61118788849522
$ws.Cells.Item(19, 3).Value = 10.02869179229833
$ws.Cells.Item(19, 5).Value = 21.58629117865512
$ws.Cells.Item(19, 6).Value = 38.1283269468263
$ws.Cells.Item(19, 7).Value = 3.611570391064483
$ws.Cells.Item(19, 10).Value = 7.640845692957991
$ws.Cells.Item(19, 13).Value = 20.05878549614288
$ws.Cells.Item(19, 15).Value = 19.53672839474854

$ws.Cells.Item(20, 2).Value = 14.83473382960551
$ws.Cells.Item(20, 3).Value = 10.17835699292673
$ws.Cells.Item(20, 5).Value = 21.61154786386557
$ws.Cells.Item(20, 6).Value = 38.15311386449196
$ws.Cells.Item(20, 7).Value = 3.610819922212818
$ws.Cells.Item(20, 10).Value = 7.632697803049511
$ws.Cells.Item(20, 13).Value = 20.16392403880463
$ws.Cells.Item(20, 15).Value = 19.50004197170744

$ws.Cells.Item(21, 2).Value = 15.56213943234277
$ws.Cells.Item(21, 3).Value = 10.6645213022439
$ws.Cells.Item(21, 5).Value = 21.70052121589602
$ws.Cells.Item(21, 6).Value = 38.25013182397885
$ws.Cells.Item(21, 7).Value = 3.608376634934577
$ws.Cells.Item(21, 10).Value = 7.606154203254986
$ws.Cells.Item(21, 13).Value = 20.51572233904529
$ws.Cells.Item(21, 15).Value = 19.38497356902453

$ws.Cells.Item(22, 2).Value = 16.01953816544246
$ws.Cells.Item(22, 3).Value = 10.96964010411682
$ws.Cells.Item(22, 5).Value = 21.76168491187714
$ws.Cells.Item(22, 6).Value = 38.32372961184984
$ws.Cells.Item(22, 7).Value = 3.606837811651103
$ws.Cells.Item(22, 10).Value = 7.589423956462459
$ws.Cells.Item(22, 13).Value = 20.74423420584205
$ws.Cells.Item(22, 15).Value = 19.31592505007606

$ws.Cells.Item(23, 2).Value = 15.77708379500144
$ws.Cells.Item(23, 3).Value = 10.80795728604335
$ws.Cells.Item(23, 5).Value = 21.72876831768111
$ws.Cells.Item(23, 6).Value = 38.2835200827925
$ws.Cells.Item(23, 7).Value = 3.607653820076878
$ws.Cells.Item(23, 10).Value = 7.598296873704611
$ws.Cells.Item(23, 13).Value = 20.62242260854863
$ws.Cells.Item(23, 15).Value = 19.35221118820358

$ws.Cells.Item(24, 2).Value = 14.82254583523933
$ws.Cells.Item(24, 3).Value = 10.17020038157188
$ws.Cells.Item(24, 5).Value = 21.61014528093603
$ws.Cells.Item(24, 6).Value = 38.15170090041162
$ws.Cells.Item(24, 7).Value = 3.610860839560052
$ws.Cells.Item(24, 10).Value = 7.633142108533328
$ws.Cells.Item(24, 13).Value = 20.15815495469608
$ws.Cells.Item(24, 15).Value = 19.50202591756599

$ws.Cells.Item(25, 2).Value = 13.71811692298821
$ws.Cells.Item(25, 3).Value = 9.429357016390544
$ws.Cells.Item(25, 5).Value = 21.49546127261518
$ws.Cells.Item(25, 6).Value = 38.05354170733307
$ws.Cells.Item(25, 7).Value = 3.614570911581587
$ws.Cells.Item(25, 10).Value = 7.673397380749217
$ws.Cells.Item(25, 13).Value = 19.65345308414039
$ws.Cells.Item(25, 15).Value = 19.68973470736064
